$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new date-like text into a scratch cell as a text formula so it
# lands as a plain string (not auto-coerced into a date serial), then copy
# only the *value* over to A4 - this keeps A4 on the default style, just
# like the existing A2/A3 "date-looking" text cells.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""2015.11.26"""
$scratch.Copy()
$ws.Range("A4").PasteSpecial(-4163) # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()

$ws.Range("B4").Value = 5
